# flare list, tasks, implementing event list function
#
# Adds an "event list" helper block in columns J:K:L mirroring the
# Date / Class / Begin columns (A / E / B) for every data row, formatted
# as yyyymmdd / plain-text / hh respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows present in the sheet (rows 8, 15, 28, 45 are blank separator rows).
$dataRows = @(2,3,4,5,6,7,9,10,11,12,13,14,16,17,18,19,20,21,22,23,24,25,26,27,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,46,47)
$blankRows = @(8,15,28,45)

# Set the number formats on the full J/L column ranges *before* writing any
# formulas, so the new cellXfs (numFmtId 164 "yyyymmdd" / 165 "hh") are
# created once and reused by every cell instead of Excel inventing one style
# per distinct precedent format it sees along the way.
$ws.Range("J2:J47").NumberFormat = "yyyymmdd"
$ws.Range("L2:L47").NumberFormat = "hh"

# Row 2 holds standalone (non-shared) formulas.
$ws.Range("J2").Formula = "=A2"
$ws.Range("L2").Formula = "=B2"

# Rows 3:47 share one formula definition each (J3:J47 / L3:L47), matching
# how Excel auto-shares a formula typed across a contiguous range.
$ws.Range("J3:J47").Formula = "=A3"
$ws.Range("L3:L47").Formula = "=B3"

# The four separator rows keep the J/L number-format styling but carry no
# value/formula.
foreach ($r in $blankRows) {
    $ws.Range("J$r").Formula = ""
    $ws.Range("L$r").Formula = ""
}

# Column K mirrors column E (Class) as a plain value with E's own style,
# for every populated data row.
foreach ($r in $dataRows) {
    $ws.Range("E$r").Copy() | Out-Null
    $ws.Range("K$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("K$r").Value = $ws.Range("E$r").Value2
}

$excel.CutCopyMode = 0

# Reflect the new used range / cursor position in the sheet view.
$ws.Range("L46").Select()
